# "Add files via upload" - new row added to the "tipos" list on the
# "datos" sheet: A4 = "aaaa" (new shared string), which also grows the
# Tabla2 ("tipos") table from A1:A3 to A1:A4.  Selections on both sheets
# were left at their new cursor positions (C17 on catalogo, A6 on datos).

$wb = $excel.ActiveWorkbook

# --- sheet "datos": add the new value and grow the table around it ---
$wsDatos = $wb.Worksheets.Item("datos")
$wsDatos.Range("A4").Value = "aaaa"

$tabla2 = $wsDatos.ListObjects.Item("Tabla2")
$tabla2.Resize($wsDatos.Range("A1:A4"))

$wsDatos.Range("A6").Select()

# --- sheet "catalogo": just a selection/cursor move, back to the active sheet ---
$wsCatalogo = $wb.Worksheets.Item("catalogo")
$wsCatalogo.Activate()
$wsCatalogo.Range("C17").Select()
